$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# 1) Update the letter date: "September 19, 2025" -> "September 21, 2025".
$dateIdx = Get-ParagraphIndexByText $d "September 19, 2025"
$pDate = $d.Paragraphs($dateIdx)
$pDate.Range.Find.Execute("September 19, 2025", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "September 21, 2025", 2)

# 2) Split the recipient address line ("3431 Copper Rd, Santa Clara CA 95051")
#    into two separate paragraphs: "3431 Copper Rd" and a new
#    "Santa Clara, CA 95051" line. Scope the Find to that single paragraph so
#    the identically worded "PROPERTY ADDRESS" entry inside the table further
#    down the letter is left untouched.
$addrIdx = Get-ParagraphIndexByText $d "3431 Copper Rd, Santa Clara CA 95051"
$pAddr = $d.Paragraphs($addrIdx)
$pAddr.Range.Find.Execute("3431 Copper Rd, Santa Clara CA 95051", $false, $false, $false, $false, $false, `
                           $true, 1, $false, "3431 Copper Rd^pSanta Clara, CA 95051", 2)

# The newly split-off paragraph/run does not inherit the surrounding Arial
# 11pt formatting automatically, so apply it explicitly to match the rest of
# the address block.
$pCity = $d.Paragraphs($addrIdx + 1)
$pCity.Range.Font.Name = "Arial"
$pCity.Range.Font.NameAscii = "Arial"
$pCity.Range.Font.NameBi = "Arial"
$pCity.Range.Font.Size = 11
$pCity.Range.Font.SizeBi = 11

# 3) Remove the now-superfluous empty "No Spacing" paragraph that immediately
#    follows "... Board of Directors" near the end of the letter.
$boardIdx = Get-ParagraphIndexByText $d "Townhomes at Nuevo Homeowners Association Board of Directors"
if ($boardIdx -gt 0) {
    $pNext = $d.Paragraphs($boardIdx + 1)
    $nextTxt = $pNext.Range.Text.TrimEnd([char]13, [char]7)
    if ($nextTxt -eq "") {
        $pNext.Range.Delete()
    }
}
